# Apply updated TPM-derived values to the Cxcl12-Itgb1 LR-pair sheet
# (recomputed Ligand/Receptor/Edge expression metrics per commit "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 119.0815153333333
$ws.Range("H2").Value = 357.244546
$ws.Range("I2").Value = 0.431812569872284
$ws.Range("J2").Value = 0.4318125698722839
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 20018.77298127835
$ws.Range("R2").Value = 180168.9568315052
$ws.Range("S2").Value = 0.1288607409808298
$ws.Range("T2").Value = 0.1288607409808298

# Row 3
$ws.Range("G3").Value = 119.0815153333333
$ws.Range("H3").Value = 357.244546
$ws.Range("I3").Value = 0.431812569872284
$ws.Range("J3").Value = 0.4318125698722839
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 19411.02955196911
$ws.Range("R3").Value = 174699.265967722
$ws.Range("S3").Value = 0.1249486995834745
$ws.Range("T3").Value = 0.1249486995834744

# Row 4
$ws.Range("G4").Value = 119.0815153333333
$ws.Range("H4").Value = 357.244546
$ws.Range("I4").Value = 0.431812569872284
$ws.Range("J4").Value = 0.4318125698722839
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 19766.76108792913
$ws.Range("R4").Value = 177900.8497913622
$ws.Range("S4").Value = 0.1272385416910265
$ws.Range("T4").Value = 0.1272385416910265

# Row 5
$ws.Range("G5").Value = 119.0815153333333
$ws.Range("H5").Value = 357.244546
$ws.Range("I5").Value = 0.431812569872284
$ws.Range("J5").Value = 0.4318125698722839
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 7886.379879991428
$ws.Range("R5").Value = 70977.41891992284
$ws.Range("S5").Value = 0.05076458761695331
$ws.Range("T5").Value = 0.05076458761695331

# Row 6
$ws.Range("I6").Value = 0.4460879372303943
$ws.Range("J6").Value = 0.4460879372303942
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 20680.57710256849
$ws.Range("R6").Value = 186125.1939231164
$ws.Range("S6").Value = 0.1331207707805267
$ws.Range("T6").Value = 0.1331207707805267

# Row 7
$ws.Range("I7").Value = 0.4460879372303943
$ws.Range("J7").Value = 0.4460879372303942
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.1290794005216149
$ws.Range("T7").Value = 0.1290794005216148

# Row 8
$ws.Range("I8").Value = 0.4460879372303943
$ws.Range("J8").Value = 0.4460879372303942
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 20420.23390390956
$ws.Range("R8").Value = 183782.1051351861
$ws.Range("S8").Value = 0.1314449429203535
$ws.Range("T8").Value = 0.1314449429203535

# Row 9
$ws.Range("I9").Value = 0.4460879372303943
$ws.Range("J9").Value = 0.4460879372303942
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 8147.097093355053
$ws.Range("R9").Value = 73323.87384019548
$ws.Range("S9").Value = 0.05244282300789926
$ws.Range("T9").Value = 0.05244282300789926

# Row 10
$ws.Range("G10").Value = 33.50679633333333
$ws.Range("H10").Value = 100.520389
$ws.Range("I10").Value = 0.1215021138451521
$ws.Range("J10").Value = 0.121502113845152
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 5632.821746089833
$ws.Range("R10").Value = 50695.3957148085
$ws.Range("S10").Value = 0.03625844524501502
$ws.Range("T10").Value = 0.03625844524501502

# Row 11
$ws.Range("G11").Value = 33.50679633333333
$ws.Range("H11").Value = 100.520389
$ws.Range("I11").Value = 0.1215021138451521
$ws.Range("J11").Value = 0.121502113845152
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 5461.816739546336
$ws.Range("R11").Value = 49156.35065591702
$ws.Range("S11").Value = 0.0351576868780944
$ws.Range("T11").Value = 0.03515768687809439

# Row 12
$ws.Range("G12").Value = 33.50679633333333
$ws.Range("H12").Value = 100.520389
$ws.Range("I12").Value = 0.1215021138451521
$ws.Range("J12").Value = 0.121502113845152
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 5561.911402361056
$ws.Range("R12").Value = 50057.20262124951
$ws.Range("S12").Value = 0.03580199571913045
$ws.Range("T12").Value = 0.03580199571913043

# Row 13
$ws.Range("G13").Value = 33.50679633333333
$ws.Range("H13").Value = 100.520389
$ws.Range("I13").Value = 0.1215021138451521
$ws.Range("J13").Value = 0.121502113845152
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 2219.045699128774
$ws.Range("R13").Value = 19971.41129215897
$ws.Range("S13").Value = 0.01428398600291222
$ws.Range("T13").Value = 0.01428398600291221

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.16474
$ws.Range("H14").Value = 0.49422
$ws.Range("I14").Value = 0.000597379052169715
$ws.Range("J14").Value = 0.000597379052169715
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 27.69441295489333
$ws.Range("R14").Value = 249.24971659404
$ws.Range("S14").Value = 0.0001782687968805147
$ws.Range("T14").Value = 0.0001782687968805147

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.16474
$ws.Range("H15").Value = 0.49422
$ws.Range("I15").Value = 0.000597379052169715
$ws.Range("J15").Value = 0.000597379052169715
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 26.85364726372667
$ws.Range("R15").Value = 241.68282537354
$ws.Range("S15").Value = 0.000172856792355746
$ws.Range("T15").Value = 0.000172856792355746

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.16474
$ws.Range("H16").Value = 0.49422
$ws.Range("I16").Value = 0.000597379052169715
$ws.Range("J16").Value = 0.000597379052169715
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 27.3457741322
$ws.Range("R16").Value = 246.1119671898
$ws.Range("S16").Value = 0.0001760246105325821
$ws.Range("T16").Value = 0.0001760246105325821

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.16474
$ws.Range("H17").Value = 0.49422
$ws.Range("I17").Value = 0.000597379052169715
$ws.Range("J17").Value = 0.000597379052169715
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 10.91019221407333
$ws.Range("R17").Value = 98.19172992666
$ws.Range("S17").Value = 0.00007022885240087239
$ws.Range("T17").Value = 0.00007022885240087239

